# Update "想去人数" (F column) counts on sheets "展览", "演出" and "全部类型"
# to match the regenerated gh-pages data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 1000
$ws1.Range("F4").Value  = 240
$ws1.Range("F6").Value  = 1164
$ws1.Range("F7").Value  = 950
$ws1.Range("F8").Value  = 286
$ws1.Range("F12").Value = 327
$ws1.Range("F14").Value = 531
$ws1.Range("F15").Value = 1380
$ws1.Range("F17").Value = 1284
$ws1.Range("F19").Value = 319
$ws1.Range("F20").Value = 1569
$ws1.Range("F21").Value = 1321
$ws1.Range("F22").Value = 759
$ws1.Range("F26").Value = 1080
$ws1.Range("F27").Value = 377
$ws1.Range("F28").Value = 3339
$ws1.Range("F29").Value = 651
$ws1.Range("F31").Value = 1476

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F15").Value = 9

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 1000
$ws4.Range("F7").Value  = 240
$ws4.Range("F10").Value = 1164
$ws4.Range("F11").Value = 950
$ws4.Range("F12").Value = 286
$ws4.Range("F24").Value = 327
$ws4.Range("F26").Value = 531
$ws4.Range("F27").Value = 1380
$ws4.Range("F29").Value = 1284
$ws4.Range("F31").Value = 319
$ws4.Range("F32").Value = 1569
$ws4.Range("F33").Value = 1321
$ws4.Range("F34").Value = 759
$ws4.Range("F40").Value = 1080
$ws4.Range("F41").Value = 377
$ws4.Range("F42").Value = 3339
$ws4.Range("F43").Value = 651
$ws4.Range("F45").Value = 1476
$ws4.Range("F47").Value = 9
